# Update weekly crime data for the 52nd Precinct CompStat report
# - Bump report volume number and covering-week date range in the header
# - Refresh all weekly/28-day/YTD/2-year crime statistics (rows 14-30)

function Set-TextCell($ws, $ref, $text) {
    # Force the cell to store the value as text (matches cells in the sheet
    # that display "0" or "***.*" placeholders instead of a number)
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 31   Number  28" -> "...Number  29" ---
$volCell = $ws.Range("A8")
$volText = $volCell.Value2
$numStart = $volText.IndexOf("28") + 1
$volChars = $volCell.Characters($numStart, 2)
$volChars.Text = "29"

# --- Header: "Report Covering the Week  7/8/2024  Through  7/14/2024" ---
$weekCell = $ws.Range("C9")
$weekText = $weekCell.Value2
$d1Start = $weekText.IndexOf("7/8/2024") + 1
$d1Chars = $weekCell.Characters($d1Start, 8)
$d1Chars.Text = "7/15/2024"

$weekCell2 = $ws.Range("C9")
$weekText2 = $weekCell2.Value2
$d2Start = $weekText2.IndexOf("7/14/2024") + 1
$d2Chars = $weekCell2.Characters($d2Start, 9)
$d2Chars.Text = "7/21/2024"

# --- Crime statistics table (rows 14-30, columns C:N) ---
    $ws.Range("L14").Value = -40
    $ws.Range("M14").Value = 0
    $ws.Range("N14").Value = -75
    Set-TextCell $ws "C15" "0"
    $ws.Range("E15").Value = -100
    $ws.Range("G15").Value = 4
    $ws.Range("H15").Value = -75
    $ws.Range("J15").Value = 19
    $ws.Range("K15").Value = 36.842105263157
    $ws.Range("L15").Value = 13.043478260869
    $ws.Range("N15").Value = -23.529411764705
    $ws.Range("D16").Value = 16
    $ws.Range("E16").Value = -6.25
    $ws.Range("F16").Value = 47
    $ws.Range("G16").Value = 51
    $ws.Range("H16").Value = -7.843137254901
    $ws.Range("I16").Value = 279
    $ws.Range("J16").Value = 254
    $ws.Range("K16").Value = 9.842519685039
    $ws.Range("L16").Value = 15.289256198347
    $ws.Range("M16").Value = 2.952029520295
    $ws.Range("N16").Value = -67.782909930715
    $ws.Range("C17").Value = 14
    $ws.Range("D17").Value = 13
    $ws.Range("E17").Value = 7.692307692307
    $ws.Range("F17").Value = 78
    $ws.Range("G17").Value = 64
    $ws.Range("H17").Value = 21.875
    $ws.Range("I17").Value = 420
    $ws.Range("J17").Value = 368
    $ws.Range("K17").Value = 14.130434782608
    $ws.Range("L17").Value = 19.318181818181
    $ws.Range("M17").Value = 69.354838709677
    $ws.Range("N17").Value = 21.037463976945
    $ws.Range("C18").Value = 6
    $ws.Range("D18").Value = 9
    $ws.Range("E18").Value = -33.333333333333
    $ws.Range("F18").Value = 15
    $ws.Range("G18").Value = 17
    $ws.Range("H18").Value = -11.764705882352
    $ws.Range("I18").Value = 107
    $ws.Range("J18").Value = 133
    $ws.Range("K18").Value = -19.548872180451
    $ws.Range("L18").Value = -10.833333333333
    $ws.Range("M18").Value = -48.557692307692
    $ws.Range("N18").Value = -91.391794046661
    $ws.Range("C19").Value = 14
    $ws.Range("D19").Value = 17
    $ws.Range("E19").Value = -17.647058823529
    $ws.Range("F19").Value = 68
    $ws.Range("G19").Value = 56
    $ws.Range("H19").Value = 21.428571428571
    $ws.Range("I19").Value = 429
    $ws.Range("J19").Value = 394
    $ws.Range("K19").Value = 8.883248730964
    $ws.Range("L19").Value = 18.181818181818
    $ws.Range("M19").Value = 49.477351916376
    $ws.Range("N19").Value = -12.269938650306
    $ws.Range("C20").Value = 10
    $ws.Range("D20").Value = 13
    $ws.Range("E20").Value = -23.076923076923
    $ws.Range("F20").Value = 27
    $ws.Range("G20").Value = 28
    $ws.Range("H20").Value = -3.571428571428
    $ws.Range("I20").Value = 193
    $ws.Range("J20").Value = 187
    $ws.Range("K20").Value = 3.208556149732
    $ws.Range("L20").Value = 28.666666666666
    $ws.Range("M20").Value = 94.949494949494
    $ws.Range("N20").Value = -77.453271028037
    $ws.Range("C21").Value = 59
    $ws.Range("D21").Value = 69
    $ws.Range("E21").Value = -14.492753623188
    $ws.Range("F21").Value = 236
    $ws.Range("G21").Value = 221
    $ws.Range("H21").Value = 6.787330316742
    $ws.Range("I21").Value = 1460
    $ws.Range("J21").Value = 1358
    $ws.Range("K21").Value = 7.511045655375
    $ws.Range("L21").Value = 15.873015873015
    $ws.Range("M21").Value = 28.747795414462
    $ws.Range("N21").Value = -62.166364343094
    Set-TextCell $ws "C22" "0"
    $ws.Range("E22").Value = -100
    $ws.Range("J22").Value = 28
    $ws.Range("K22").Value = -10.714285714285
    $ws.Range("L22").Value = 4.166666666666
    $ws.Range("C24").Value = 26
    $ws.Range("D24").Value = 30
    $ws.Range("E24").Value = -13.333333333333
    $ws.Range("F24").Value = 122
    $ws.Range("G24").Value = 154
    $ws.Range("H24").Value = -20.779220779220
    $ws.Range("I24").Value = 1012
    $ws.Range("J24").Value = 1419
    $ws.Range("K24").Value = -28.682170542635
    $ws.Range("L24").Value = -38.405356055995
    $ws.Range("M24").Value = 41.340782122905
    $ws.Range("C25").Value = 12
    $ws.Range("D25").Value = 18
    $ws.Range("E25").Value = -33.333333333333
    $ws.Range("F25").Value = 54
    $ws.Range("G25").Value = 96
    $ws.Range("H25").Value = -43.75
    $ws.Range("I25").Value = 554
    $ws.Range("J25").Value = 943
    $ws.Range("K25").Value = -41.251325556733
    $ws.Range("L25").Value = -55.502008032128
    $ws.Range("C26").Value = 14
    $ws.Range("D26").Value = 16
    $ws.Range("E26").Value = -12.5
    $ws.Range("F26").Value = 79
    $ws.Range("G26").Value = 61
    $ws.Range("H26").Value = 29.508196721311
    $ws.Range("I26").Value = 479
    $ws.Range("J26").Value = 461
    $ws.Range("K26").Value = 3.904555314533
    $ws.Range("L26").Value = 6.681514476614
    $ws.Range("M26").Value = -3.815261044176
    Set-TextCell $ws "C27" "0"
    $ws.Range("E27").Value = -100
    $ws.Range("F27").Value = 2
    $ws.Range("G27").Value = 4
    $ws.Range("H27").Value = -50
    $ws.Range("J27").Value = 29
    $ws.Range("K27").Value = 31.034482758620
    $ws.Range("L27").Value = -22.448979591836
    $ws.Range("C28").Value = 2
    Set-TextCell $ws "D28" "0"
    Set-TextCell $ws "E28" "***.*"
    $ws.Range("F28").Value = 6
    $ws.Range("G28").Value = 8
    $ws.Range("H28").Value = -25
    $ws.Range("I28").Value = 66
    $ws.Range("K28").Value = 20
    $ws.Range("L28").Value = 60.975609756097
    Set-TextCell $ws "D29" "0"
    Set-TextCell $ws "E29" "***.*"
    $ws.Range("F29").Value = 3
    $ws.Range("H29").Value = -40
    $ws.Range("L29").Value = -13.636363636363
    $ws.Range("M29").Value = 11.764705882352
    $ws.Range("N29").Value = -53.658536585365
    Set-TextCell $ws "D30" "0"
    Set-TextCell $ws "E30" "***.*"
    $ws.Range("F30").Value = 2
    $ws.Range("H30").Value = 0
    $ws.Range("L30").Value = -36.842105263157
    $ws.Range("M30").Value = -20
    $ws.Range("N30").Value = -66.666666666666
